$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename TestBean to JavaBean to avoid interpolation with JUnit tests in surefire plugin
$ws.Range("F5").Value = "Data JavaBean beans"
$ws.Range("B3").Value = "Method String print(JavaBean bean)"

$ws.Range("B4").Select()
